$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.26%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'41.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'4.27%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.113"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.30%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07621"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-0.67%"
$ws.Range("E5").Style = "Normal"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.248"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.15%"
$ws.Range("E6").Style = "Normal"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.606"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.46%"
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.460"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.67%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9019"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.01%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1121"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'10.94%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1792"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'3.55%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09229"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'3.67%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.04160"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-6.76%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.1052"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.17%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001255"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-1.11%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005873"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.37%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.340"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.42%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.3323"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.24%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'6.627"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-6.57%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'1.04%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.2704"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.04068"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-2.13%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.001245"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'2.92%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004090"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.00%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'-0.22%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0003748"
$ws.Range("D26").Style = "Normal"
$ws.Range("D38").Value = "'0.02400"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'3.46%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05225"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'1.50%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007781"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.71%"
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'-1.47%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.007057"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'8.88%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'-0.60%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.007708"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-6.50%"
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'-7.46%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006944"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'6.28%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.32%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.05383"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'1,454.50%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.32%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.32%"
$ws.Range("E51").Style = "Normal"
